# Volt Var debugging complete
# Update the PFlim (column N) setpoint on the "Controllers" sheet from
# 0.85 to 0.9 for every data row (rows 3 through 533), and leave the
# selection on that column/range, matching the author's last action.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Controllers")

$firstRow = 3
$lastRow = 533
$col = 14   # column N = PFlim

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, $col).Value = 0.9
}

$ws.Range("N3:N533").Select()
